# Continuação alteração dos dados, script python para inserção e SQL Programming
#
# - Rename the "Marcação" sheet to "marcacao"
# - Rename the header of its last column (E1) from "data" to "dataMarc"
# - Make that sheet the active tab, with E19 selected (losing the tabSelected
#   flag on the previously active "produto" sheet automatically)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Marcação")
$ws.Name = "marcacao"

$ws.Range("E1").Value = "dataMarc"

$ws.Activate()
$ws.Range("E19").Select()
